# Auto-generated: updates column F ("想去人数") values across sheets
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5030
$ws.Range("F6").Value = 5030
$ws.Range("F8").Value = 143
$ws.Range("F10").Value = 40
$ws.Range("F11").Value = 210
$ws.Range("F12").Value = 163
$ws.Range("F13").Value = 8279
$ws.Range("F14").Value = 270
$ws.Range("F15").Value = 138
$ws.Range("F17").Value = 602
$ws.Range("F18").Value = 2506
$ws.Range("F19").Value = 6323
$ws.Range("F20").Value = 2295
$ws.Range("F23").Value = 2514
$ws.Range("F24").Value = 19
$ws.Range("F25").Value = 12
$ws.Range("F26").Value = 6373
$ws.Range("F27").Value = 179
$ws.Range("F28").Value = 64
$ws.Range("F30").Value = 100
$ws.Range("F32").Value = 6829
$ws.Range("F36").Value = 9
$ws.Range("F41").Value = 43
$ws.Range("F42").Value = 2515
$ws.Range("F46").Value = 60
$ws.Range("F47").Value = 506
$ws.Range("F48").Value = 2206
$ws.Range("F49").Value = 66
$ws.Range("F50").Value = 1117

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 155
$ws.Range("F8").Value = 12
$ws.Range("F14").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5030
$ws.Range("F5").Value = 5030
$ws.Range("F7").Value = 143
$ws.Range("F9").Value = 40
$ws.Range("F10").Value = 210
$ws.Range("F11").Value = 163
$ws.Range("F12").Value = 8279
$ws.Range("F13").Value = 8279
$ws.Range("F14").Value = 270
$ws.Range("F15").Value = 138
$ws.Range("F16").Value = 602
$ws.Range("F17").Value = 2506
$ws.Range("F19").Value = 155
$ws.Range("F20").Value = 6323
$ws.Range("F21").Value = 2295
$ws.Range("F24").Value = 2514
$ws.Range("F25").Value = 19
$ws.Range("F27").Value = 12
$ws.Range("F28").Value = 6373
$ws.Range("F29").Value = 179
$ws.Range("F30").Value = 64
$ws.Range("F32").Value = 100
$ws.Range("F34").Value = 6829
$ws.Range("F41").Value = 2515
$ws.Range("F44").Value = 60
$ws.Range("F45").Value = 506
$ws.Range("F47").Value = 2206
$ws.Range("F48").Value = 66
$ws.Range("F50").Value = 1117
$ws.Range("F51").Value = 26
